$d = $word.ActiveDocument

$replacements = @(
    @{old="42×82=3444"; new="33×37=1221"},
    @{old="35×73=2555"; new="54×25=1350"},
    @{old="90×52=4680"; new="14×49=686"},
    @{old="36×74=2664"; new="35×18=630"},
    @{old="23×62=1426"; new="44×98=4312"},
    @{old="81×25=2025"; new="74×34=2516"},
    @{old="95×67=6365"; new="82×61=5002"},
    @{old="44×84=3696"; new="13×98=1274"},
    @{old="74×79=5846"; new="96×72=6912"},
    @{old="56×55=3080"; new="46×86=3956"},
    @{old="59×91=5369"; new="75×34=2550"},
    @{old="87×91=7917"; new="87×20=1740"},
    @{old="94×18=1692"; new="14×15=210"},
    @{old="26×23=598"; new="49×27=1323"},
    @{old="13×12=156"; new="28×40=1120"},
    @{old="96×26=2496"; new="90×53=4770"},
    @{old="30×93=2790"; new="67×11=737"},
    @{old="72×97=6984"; new="55×64=3520"},
    @{old="12×49=588"; new="63×93=5859"},
    @{old="55×37=2035"; new="48×56=2688"},
    @{old="38×93=3534"; new="92×24=2208"},
    @{old="14×94=1316"; new="60×58=3480"},
    @{old="29×73=2117"; new="60×29=1740"},
    @{old="58×70=4060"; new="98×36=3528"},
    @{old="89×37=3293"; new="79×57=4503"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false,
                         $true, 1, $false, $r.new, 2)
}

$d.Save()
